$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert three new rows before the existing data row (row 3), pushing the
# original data row ("Câmera DSLR...") down to row 6.
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()

# Rows 3, 4 and 5 repeat the header labels (same text as row 2), but will
# use the data-row style (copied below from row 6, the former row 3).
$headers = @("Nº DE ORDEM","TOMBO","DESCRIÇÃO DO BEM","DATA DA AQUISIÇÃO","DOCUMENTO FISCAL","UNIDADE RESPONSÁVEL","CLASSIFICAÇÃO","DESTINAÇÃO")

foreach ($r in 3,4,5) {
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value = $headers[$c - 1]
    }
}

# New data rows that go after the original one (rows 7-11).
$data = @(
    @(3, 100103, "Mesa de escritório 1.20m x 0.60m, cor carvalho", "20/11/2021", "NF-e 33210", "Diretoria de Comunicação (DECOM)", "Irrecuperável", "Alienação/Leilão"),
    @(4, 199887, "Maca hospitalar simples com rodízios", "15/12/2021", "NF-e 34567", "Centro de Ciências da Saúde e do Desporto (CCSD)", "Irrecuperável", "Alienação/Leilão"),
    @(5, 134567, "Arquivo de aço com 4 gavetas, para pasta suspensa", "19/05/2017", "NF-e 8123", "Diretoria de Sistemas de Informação (DSI)", "Irrecuperável", "Alienação/Leilão"),
    @(7, 112233, "Multímetro Digital Minipa ET-2042E", "10/05/2023", "NF-e 54001", "Coordenação do Curso de Engenharia Elétrica", "Irrecuperável", "Alienação/Leilão"),
    @(9, 121212, "Switch de rede 24 portas Gigabit, TP-Link", "07/07/2022", "NF-e 48500", "Pró-Reitoria de Desenvolvimento e Gestão de Pessoas (PRODGEP)", "Irrecuperável", "Alienação/Leilão")
)

# Force the "DATA DA AQUISIÇÃO" column to be stored as plain text so values
# such as "10/05/2023" are kept as text instead of being auto-converted
# into date serial numbers.
$ws.Range("D7:D11").NumberFormat = "@"

$r = 7
foreach ($row in $data) {
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r = $r + 1
}

# Copy the formatting (style, incl. borders/alignment) of the original data
# row (now row 6) onto the new header-like rows (3-5) and the new data rows
# (7-11), and match their row height too.
$fmtSrc = $ws.Range("A6:H6")
$fmtSrc.Copy()
foreach ($r in 3,4,5,7,8,9,10,11) {
    $dst = $ws.Range($ws.Cells.Item($r,1), $ws.Cells.Item($r,8))
    $dst.PasteSpecial(-4122)
    $ws.Rows.Item($r).RowHeight = $ws.Rows.Item(6).RowHeight
}
$excel.CutCopyMode = 0

[void]$ws.Range("A1").Select()
